$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 4 (old PINT2/SICA/PINTURA2 -> new 15569/Chiesa/Carbon G720)
# Force A4 to remain text (not be reinterpreted as a number) like the rest of the "Codigo" column
$ws.Range("A4").NumberFormat = "@"
$ws.Range("A4").Value = "15569"
$ws.Range("A4").Style = "Normal"
$ws.Range("B4").Value = "Chiesa"
$ws.Range("C4").Value = "Carbon G720"

# Update row 5 (old AMOLA1/DEWALT/AMOLAD1/100/ricky -> new PINT6/Sorbalok/Pintura33/0/ferreteria_general)
$ws.Range("A5").Value = "PINT6"
$ws.Range("B5").Value = "Sorbalok"
$ws.Range("C5").Value = "Pintura33"
$ws.Range("D5").Value = 0
$ws.Range("F5").Value = "ferreteria_general"

# Add (empty) "Observaciones" cells in column E for rows 2-5
$ws.Range("E2").Value = ""
$ws.Range("E3").Value = ""
$ws.Range("E4").Value = ""
$ws.Range("E5").Value = ""

# Remove old rows 6 and 7 (FOC1/... and AMOLA2/... entries), shrinking the sheet to A1:F5
$ws.Rows("6:7").Delete()
